{"js": "// Update the scorecardBatch.scorecardUrl value in the document:\n// old: https://devccda.sitenv.org/scorecard/ccdascorecardservice2\n// new: https://ccda.healthit.gov/scorecard/ccdascorecardservice2\nconst oldUrl = \"https://devccda.sitenv.org/scorecard/ccdascorecardservice2\";\nconst newUrl = \"https://ccda.healthit.gov/scorecard/ccdascorecardservice2\";\n\nconst results = context.document.body.search(oldUrl, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newUrl, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the scorecardBatch.scorecardUrl value in the document:\n#   old: https://devccda.sitenv.org/scorecard/ccdascorecardservice2\n#   new: https://ccda.healthit.gov/scorecard/ccdascorecardservice2\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Execute(\n  \"https://devccda.sitenv.org/scorecard/ccdascorecardservice2\",\n  $false,\n  $false,\n  $false,\n  $false,\n  $false,\n  $true,\n  1,\n  $false,\n  \"https://ccda.healthit.gov/scorecard/ccdascorecardservice2\",\n  2\n)\n"}
